$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EXE Logout")
$ws.Activate()

# Copy the formatting of row 6 (an existing "logout" test case row) down to the
# new row 7 so the new row picks up the same cell styles (borders/fills/fonts).
$ws.Range("D6:Q6").Copy() | Out-Null
$ws.Range("D7:Q7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new test case ("Verify user session is terminated after logout").
$ws.Range("D7").Value = "EXE-03"
$ws.Range("E7").Value = "Verify user session is terminated after logout"
$ws.Range("F7").Value = "EXE-TC-03"
$ws.Range("G7").Value = "Logout Module"
$ws.Range("H7").Value = "EXE-03"
$ws.Range("I7").Value = "Verify session is invalid after logout"
$ws.Range("J7").Value = "User cannot access restricted page"
$ws.Range("K7").Value = "User cannot access restricted paage after logout"
$ws.Range("L7").Value = "PASS"
$ws.Range("M7").Value = "-"
$ws.Range("N7").Value = "-"
$ws.Range("O7").Value = "Syaif (QA)"
$ws.Range("P7").Value = 46081
$ws.Range("Q7").Value = "Chrome v145 /`nWindows 13"

# Match the row's visual height to the other wrapped-text rows on this sheet.
$ws.Rows.Item(7).RowHeight = 31.5

# Move the view back to the top-left of the sheet and select the new cell
# that was just typed in (K7), matching the saved selection state.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K7").Select() | Out-Null
